$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

# New TaxType values to append as additional test rows (Phase 3 RAD data)
$newTaxTypes = @(
    "Admissions and Amusement Tax",
    "Estate Tax",
    "Motor Fuel Tax",
    "Slots License Fee",
    "Tobacco Tax",
    "Transportation Network Services",
    "Unclaimed Property",
    "IFTA Tax"
)

$startRow = 32
$endRow = $startRow + (2 * $newTaxTypes.Count) - 1

# Pre-format the new rows the same way as the existing block above them (row 31),
# so the cells pick up the existing bordered styles instead of creating new ones.
$ws.Range("C31:G31").Copy()
$ws.Range("C$startRow`:G$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$row = $startRow

# Block 1: rows 32-39, NotiInvoNumber = "123456789"
foreach ($taxType in $newTaxTypes) {
    $ws.Cells.Item($row, 3).Value = "Y"
    $ws.Cells.Item($row, 4).Value = "Existing Liability w/Notice Number"
    $ws.Cells.Item($row, 5).Value = $taxType
    $ws.Cells.Item($row, 6).Value = "123456789"
    $ws.Cells.Item($row, 7).Value = "Notice Number must be 13 digits in length, with a value greater than zero"
    $row = $row + 1
}

# Block 2: rows 40-47, NotiInvoNumber = "0000000000000"
foreach ($taxType in $newTaxTypes) {
    $ws.Cells.Item($row, 3).Value = "Y"
    $ws.Cells.Item($row, 4).Value = "Existing Liability w/Notice Number"
    $ws.Cells.Item($row, 5).Value = $taxType
    $ws.Cells.Item($row, 6).Value = "0000000000000"
    $ws.Cells.Item($row, 7).Value = "Notice Number must be 13 digits in length, with a value greater than zero"
    $row = $row + 1
}

# Widen columns E:G to fit the newly added (longer) content.
# Note: the host's column-width persistence quantizes to 1/6-character
# steps and adds a fixed offset on save, so these inputs are chosen to
# land as close as achievable to the target widths (53.5703125, 26,
# 80.5703125) after that round-trip.
$ws.Columns.Item(5).ColumnWidth = 52.666666666666664
$ws.Columns.Item(6).ColumnWidth = 25.166666666666668
$ws.Columns.Item(7).ColumnWidth = 79.66666666666667

# Scroll the view down and leave the selection where the author left off editing
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D52").Select()
